$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 773, shifting existing rows 773:866 down to 774:867
$ws.Rows("773:773").Insert()

# Populate the newly inserted row 773 with the new record's data
$ws.Cells.Item(773, 1).Value = 10
$ws.Cells.Item(773, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(773, 3).Value = "La Araucanía"
$ws.Cells.Item(773, 4).Value = 45124
$ws.Cells.Item(773, 5).Value = 9
$ws.Cells.Item(773, 6).Value = 100112045
$ws.Cells.Item(773, 7).Value = "Zapallo"
$ws.Cells.Item(773, 8).Value = "Camote"
$ws.Cells.Item(773, 9).Value = "1a (guarda)"
$ws.Cells.Item(773, 10).Value = 1200
$ws.Cells.Item(773, 11).Value = 600
$ws.Cells.Item(773, 12).Value = 600
$ws.Cells.Item(773, 13).Value = 600
$ws.Cells.Item(773, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(773, 15).Value = "Región del Maule"
$ws.Cells.Item(773, 16).Value = 600
$ws.Cells.Item(773, 17).Value = 1
$ws.Cells.Item(773, 18).Value = "Hortaliza"

# Apply the same date number format as the other date cells in column D
$ws.Cells.Item(773, 4).NumberFormat = $ws.Cells.Item(774, 4).NumberFormat
